$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1, formatted to match the existing header row (E1)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# New column F values: time_taken timestamps for each data row
$ws.Range("F2").Value = "2021-10-05 13:38:50.791359"
$ws.Range("F3").Value = "2021-10-05 13:38:50.791366"
$ws.Range("F4").Value = "2021-10-05 13:38:50.791369"
$ws.Range("F5").Value = "2021-10-05 13:38:50.791371"
$ws.Range("F6").Value = "2021-10-05 13:38:50.791373"
